$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: add new entry "543. Diameter of Binary Tree"
$ws.Range("A8").Value = 543
$ws.Range("B8").Value = "LC"
$ws.Range("C8").Value = "Diameter of Binary Tree"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "Java/Python"
$ws.Range("D8").VerticalAlignment = -4160

# Update selection to D8
$ws.Range("D8").Select()
